# Fruta / hortaliza, semanal
# Insert a new data row at row 47 (pushing the existing rows 47-104 down to
# 48-105) and populate the new row with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 47:104 down by one, creating a blank row 47.
$ws.Rows("47:47").Insert()

# Fill in the new row 47 with the new weekly record.
$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = 44771
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 100112035
$ws.Range("G47").Value = "Bruselas (repollito)"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 80
$ws.Range("K47").Value = 25000
$ws.Range("L47").Value = 26000
$ws.Range("M47").Value = 25500
$ws.Range("N47").Value = "`$/malla 10 kilos"
$ws.Range("O47").Value = "Región Metropolitana"
$ws.Range("P47").Value = 2550
$ws.Range("Q47").Value = 10
$ws.Range("R47").Value = "Hortaliza"
